$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,
                                            estimator=Pipeline(steps=[('model',
                                                                       GradientBoostingRegressor())]),
                                            param_grid={'model__max_depth': [3,
                                                                             5,
                                                                             7],
                                                        'model__n_estimators': [50,
                                                                                100,
                                                                                150]},
                                            scoring='neg_mean_squared_error'))"

# New header cell F1 - copy header style from A1, then set the text
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# New data cells F2/F3 with model description
$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText

# Minor precision corrections in B3/D3
$ws.Range("B3").Value = 0.08311721650820843
$ws.Range("D3").Value = 0.2072673588334912
